$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 2; D = '26.052.95'; E = '  +0.60%  ' },
    @{ Row = 3; D = '1.644.03'; E = '  +0.30%  ' },
    @{ Row = 4; D = $null; E = '  -0.29%  ' },
    @{ Row = 5; D = '214.74'; E = '  -0.10%  ' },
    @{ Row = 6; D = '0.5098'; E = '  +1.52%  ' },
    @{ Row = 7; D = $null; E = '  -0.26%  ' },
    @{ Row = 8; D = '0.2564'; E = '  -0.01%  ' },
    @{ Row = 9; D = '0.06357'; E = '  -0.51%  ' },
    @{ Row = 10; D = $null; E = '  -0.05%  ' },
    @{ Row = 11; D = '0.07758'; E = '  -0.40%  ' },
    @{ Row = 12; D = '4.291'; E = '  +0.40%  ' },
    @{ Row = 13; D = '1.641.11'; E = '  -0.58%  ' },
    @{ Row = 14; D = '0.5447'; E = '  +0.35%  ' },
    @{ Row = 15; D = '64.33'; E = '  -0.51%  ' },
    @{ Row = 16; D = '0.0₅7728'; E = '  -1.64%  ' },
    @{ Row = 17; D = '26.061.08'; E = '  +0.59%  ' },
    @{ Row = 18; D = $null; E = '  -0.31%  ' },
    @{ Row = 19; D = '198.96'; E = '  +0.46%  ' },
    @{ Row = 20; D = '4.428'; E = '  +0.98%  ' },
    @{ Row = 21; D = $null; E = '  -0.06%  ' },
    @{ Row = 22; D = '6.043'; E = '  +1.21%  ' },
    @{ Row = 23; D = $null; E = '  -0.28%  ' },
    @{ Row = 24; D = '1.869'; E = '  -0.16%  ' },
    @{ Row = 25; D = '140.85'; E = '  +0.69%  ' },
    @{ Row = 26; D = '0.1198'; E = '  +5.06%  ' },
    @{ Row = 27; D = '6.818'; E = '  -0.26%  ' },
    @{ Row = 28; D = '15.57'; E = '  -0.56%  ' },
    @{ Row = 29; D = '1.235'; E = '  -0.51%  ' },
    @{ Row = 30; D = '0.04862'; E = '  -0.50%  ' },
    @{ Row = 31; D = '3.258'; E = '  +0.07%  ' },
    @{ Row = 32; D = '3.168'; E = '  -0.70%  ' },
    @{ Row = 33; D = $null; E = '  -0.21%  ' },
    @{ Row = 34; D = '2.363'; E = '  -0.21%  ' },
    @{ Row = 35; D = '0.8996'; E = '  +1.00%  ' },
    @{ Row = 36; D = $null; E = '  -0.78%  ' },
    @{ Row = 37; D = '1.141.13'; E = '  +0.47%  ' },
    @{ Row = 38; D = '0.5474'; E = '  -1.24%  ' },
    @{ Row = 39; D = '0.01567'; E = '  +0.42%  ' },
    @{ Row = 40; D = '1.001'; E = '  -0.45%  ' },
    @{ Row = 41; D = $null; E = '  -0.71%  ' },
    @{ Row = 42; D = $null; E = '  +4.36%  ' },
    @{ Row = 43; D = '0.8117'; E = '  -0.54%  ' },
    @{ Row = 44; D = '99.41'; E = '  -0.11%  ' },
    @{ Row = 45; D = '5.389'; E = '  -5.23%  ' },
    @{ Row = 46; D = '1.780.90'; E = '  +0.42%  ' },
    @{ Row = 47; D = $null; E = '  +0.11%  ' },
    @{ Row = 48; D = '54.97'; E = '  -0.63%  ' },
    @{ Row = 49; D = '1.001'; E = '  -0.47%  ' },
    @{ Row = 50; D = '0.05054'; E = '  -0.53%  ' },
    @{ Row = 51; D = '1.002'; E = '  -0.38%  ' }
)

foreach ($change in $changes) {
    $row = $change.Row
    if ($null -ne $change.D) {
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value2 = $change.D
        $dCell.Style = "Normal"
    }
    if ($null -ne $change.E) {
        $ws.Cells.Item($row, 5).Value2 = $change.E
    }
}

"Applied $($changes.Count) row updates"
